$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.985.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.64%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.601.53'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.37%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '523.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.70%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.91%  '

# Row 7
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.594'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.94%  '

# Row 9
$ws.Range("E9").Value = '  +2.35%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.106'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.18%  '

# Row 11
$ws.Range("E11").Value = '  +0.03%  '

# Row 12
$ws.Range("E12").Value = '  +1.25%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.056.92'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.21%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '61.011.21'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.52%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.75'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.40%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000142'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.67%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.602.93'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.20%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.77'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.16%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '355.57'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.20%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.56'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.81%  '

# Row 21
$ws.Range("E21").Value = '  +2.38%  '

# Row 22
$ws.Range("E22").Value = '  +0.13%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.19%  '

# Row 24
$ws.Range("E24").Value = '  +2.21%  '

# Row 25
$ws.Range("E25").Value = '  +0.77%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.715.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.81%  '

# Row 27
$ws.Range("E27").Value = '  -0.60%  '

# Row 28
$ws.Range("D28").Value = '0.0₃0848'
$ws.Range("E28").Value = '  +1.48%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.39'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.74%  '

# Row 30
$ws.Range("E30").Value = '  -0.05%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.41'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.84%  '

# Row 33
$ws.Range("E33").Value = '  +3.23%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.12'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.30%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.20'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.92%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.932'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.90%  '

# Row 37
$ws.Range("E37").Value = '  +2.08%  '

# Row 38
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.872'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.08%  '

# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.51'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.03%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.82'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.46%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.47'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.47%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '290.51'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.62%  '

# Row 43
$ws.Range("E43").Value = '  +2.58%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.622'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.84%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0562'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.51%  '

# Row 46
$ws.Range("E46").Value = '  -0.15%  '

# Row 47
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.05'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.46%  '

# Row 48
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.15%  '

# Row 49
$ws.Range("E49").Value = '  +2.57%  '

# Row 50
$ws.Range("E50").Value = '  +0.40%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.23'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.71%  '
